$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# TC1 (CasesTab) query: drop the trailing OPTIONAL cohort-match clause and its
# `Cohort` column from the RETURN list (timing issue fix).
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and demo.sex in ['Female'] and demo.neutered_indicator in ['Yes'] OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Removing the trailing cohort line shortens the wrapped text, so the
# auto-fit row height shrinks from 304.5 to 290 (same as rows 3/4).
$ws.Rows.Item(2).RowHeight = 290

# Move the view to where the just-edited cell lives.
[void]$ws.Application.Goto($ws.Range("A2"), $true)
$ws.Range("C2").Select() | Out-Null
